$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1 - "Label", styled like the other header cells (same style as G1)
$ws.Range("H1").Value = "Label"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122) # xlPasteFormats

# Fill the new Label column (H) with 0/1 values per patient group, for both
# per-patient blocks of rows (2-11 = 100 iterations, 12-21 = 200 iterations)
$labels = @(0,0,0,0,0,1,1,1,1,1)
for ($i = 0; $i -lt 10; $i++) {
    $row2 = 2 + $i
    $row12 = 12 + $i
    $ws.Cells.Item($row2, 8).Value = $labels[$i]
    $ws.Cells.Item($row12, 8).Value = $labels[$i]
}
